$d = $word.ActiveDocument

# --- 1) TECHNICAL SKILLS / Programming Languages line -----------------
# Before: " Java, C++/C, Python, SQL, HTML, JavaScript, Bash" + ", Processing"
# After : " Java, C++/C, Python, Assembly, SQL, HTML, JavaScript, Processing "
$d.Content.Find.Execute(
    "Java, C++/C, Python, SQL, HTML, JavaScript, Bash, Processing",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Java, C++/C, Python, Assembly, SQL, HTML, JavaScript, Processing ",
    2) | Out-Null

# --- 2) Smart cane bullet: merge runs / drop proofErr spell markers ---
$d.Content.Find.Execute(
    "Designed and implemented a smart cane for the blind, by detecting objects and warning with raspberry pi4B, Google CoCo model, cam, ultrasonic sensor",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Designed and implemented a smart cane for the blind, by detecting objects and warning with raspberry pi4B, Google CoCo model, cam, ultrasonic sensor",
    2) | Out-Null

# --- 3) PERSONAL INFO line: merge the leading-space run with the text run
$d.Content.Find.Execute(
    "PERSONAL INFO -- U.S. Citizen, born in Seattle, WA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PERSONAL INFO -- U.S. Citizen, born in Seattle, WA",
    2) | Out-Null

# --- 4) Turn on "Different First Page" header/footer, which mints
#        header2/header3/footer1/footer2/footer3 parts (default + first +
#        even footers) exactly like Word does when the option is toggled
#        on and the section is touched.
$sec = $d.Sections.First
$sec.PageSetup.DifferentFirstPageHeaderFooter = $true
$sec.Footers(1).Range.Text = ""
